# Added QOL features and serialisation features
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old sample data area (columns C:J are no longer used) before
# writing the new, fuller list of store codes.
$ws.Range("A2:K12").ClearContents()

# New store-code rows: Store Code (A), Time (B), Notes (K).
# Columns C-J (Customers / Re-entry Customers / Suspected Staff / etc.)
# are intentionally left blank for this data set.
$rows = @(
    @("AU045",    11),
    @("KR045",    19),
    @("KR054",    19),
    @("KR055-01", 14),
    @("KR055-02", 14),
    @("KR057",    15),
    @("KR062",    14),
    @("MO002-02", 11),
    @("SG070",    17),
    @("TH070",     0),
    @("VN002",    21)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 11).Value = ""
    $r = $r + 1
}
